$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap content between paired rows (4<->5, 22<->23, 27<->28) ---
# (Row B holds the "Taxonsorteringsordning" sort key and a handful of
#  record fields; pairs of rows had their records swapped.)

# swap rows 4 and 5
$ws.Range("A4").Value = 130864689
$ws.Range("A5").Value = 130864687
$ws.Range("B4").Value = 78251
$ws.Range("B5").Value = 83219
$ws.Range("E4").Value = 228579
$ws.Range("E5").Value = 6440
$ws.Range("F4").Value = "Liten svartspik"
$ws.Range("F5").Value = "Vitgrynig nållav"
$ws.Range("G4").Value = "Chaenothecopsis nana"
$ws.Range("G5").Value = "Chaenotheca subroscida"
$ws.Range("H4").Value = "Tibell"
$ws.Range("H5").Value = "(Eitner) Zahlbr."
$ws.Range("Q4").Value = 446026
$ws.Range("Q5").Value = 445985
$ws.Range("R4").Value = 7031030
$ws.Range("R5").Value = 7030968
$ws.Range("S4").Value = 4
$ws.Range("S5").Value = 5
$ws.Range("Z4").Value = "13:51"
$ws.Range("Z5").Value = "11:39"
$ws.Range("AB4").Value = "13:51"
$ws.Range("AB5").Value = "11:39"

# swap rows 22 and 23
$ws.Range("A22").Value = 130864515
$ws.Range("A23").Value = 130864521
$ws.Range("B22").Value = 83219
$ws.Range("B23").Value = 91767
$ws.Range("D22").Value = "NT"
$ws.Range("D23").Value = "LC"
$ws.Range("E22").Value = 6440
$ws.Range("E23").Value = 5447
$ws.Range("F22").Value = "Vitgrynig nållav"
$ws.Range("F23").Value = "Vedticka"
$ws.Range("G22").Value = "Chaenotheca subroscida"
$ws.Range("G23").Value = "Fuscoporia viticola"
$ws.Range("H22").Value = "(Eitner) Zahlbr."
$ws.Range("H23").Value = "(Schwein.) Murrill"
$ws.Range("Q22").Value = 445932
$ws.Range("Q23").Value = 446069
$ws.Range("R22").Value = 7031103
$ws.Range("R23").Value = 7030939
$ws.Range("S22").Value = 3
$ws.Range("S23").Value = 8
$ws.Range("Z22").Value = "11:08"
$ws.Range("Z23").Value = "13:20"
$ws.Range("AB22").Value = "11:08"
$ws.Range("AB23").Value = "13:20"
$ws.Range("AC22").Value = "På bark på stam av levande gammal gran"
$ws.Range("AC23").Value = "På död klen gran i gammal granskog"

# swap rows 27 and 28
$ws.Range("A27").Value = 130864518
$ws.Range("A28").Value = 130864530
$ws.Range("B27").Value = 83219
$ws.Range("B28").Value = 79239
$ws.Range("E27").Value = 6440
$ws.Range("E28").Value = 6425
$ws.Range("F27").Value = "Vitgrynig nållav"
$ws.Range("F28").Value = "Garnlav"
$ws.Range("G27").Value = "Chaenotheca subroscida"
$ws.Range("G28").Value = "Alectoria sarmentosa"
$ws.Range("H27").Value = "(Eitner) Zahlbr."
$ws.Range("H28").Value = "(Ach.) Ach."
$ws.Range("Q27").Value = 446031
$ws.Range("Q28").Value = 445967
$ws.Range("R27").Value = 7030902
$ws.Range("R28").Value = 7031099
$ws.Range("S27").Value = 5
$ws.Range("S28").Value = 6
$ws.Range("Z27").Value = "11:56"
$ws.Range("Z28").Value = "14:21"
$ws.Range("AB27").Value = "11:56"
$ws.Range("AB28").Value = "14:21"
$ws.Range("AC27").Value = "På bark på stam av levande gammal gran i gammal granskog"
$ws.Range("AC28").Value = "På gammal gran (ca 150 år)i gammal granskog"

# Rows 22/23 also swap a few always-blank placeholder cells (K, N, AF):
# row 22 picks up row 23's (present-but-empty) blanks, row 23 reverts to
# having no cell there at all (which is what row 22 had before the swap).
$ws.Range("K22").Value = ""
$ws.Range("N22").Value = ""
$ws.Range("AF22").Value = ""
$ws.Range("K23").ClearContents()
$ws.Range("N23").ClearContents()
$ws.Range("AF23").ClearContents()

# --- Step 2: bump every "Taxonsorteringsordning" (column B) value by 4 ---
for ($r = 2; $r -le 30; $r++) {
    $cell = $ws.Range("B$r")
    $cell.Value = $cell.Value2 + 4
}
